# Apply the "Run through with latest mex files" update.
#
# The re-coded segment removed five single-letter "sub-code" rows that had
# become redundant duplicates of the "Drug Resisted" rows immediately above
# them (rows 171 "C", 173 "D", 179 "F", 181 "G", 183 "H" in the original
# layout), and appended one freshly coded row at the bottom of the table
# (a new "Drug Resisted" / teicoplanin segment coded by "chen").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the five now-redundant rows. Must go from bottom to top so that
# earlier row numbers stay valid as later deletions happen.
$rowsToDelete = @(171, 173, 179, 181, 183) | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# After the five deletions the table now runs A1:M182; append the new
# segment as row 183.
$newRow = 183
$ws.Cells.Item($newRow, 1).Value = "$([char]0x25CF)"
$ws.Cells.Item($newRow, 2).Value = ""
$ws.Cells.Item($newRow, 3).Value = ""
$ws.Cells.Item($newRow, 4).Value = "13804"
$ws.Cells.Item($newRow, 5).Value = "Drug Resisted"
$ws.Cells.Item($newRow, 6).Value = "1: 4421"
$ws.Cells.Item($newRow, 7).Value = "1: 4431"
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = "teicoplanin"
$ws.Cells.Item($newRow, 10).Value = 11
$ws.Cells.Item($newRow, 11).Value = 0.12404149751917004
$ws.Cells.Item($newRow, 12).Value = "chen"
$ws.Cells.Item($newRow, 13).Value = "1/29/19 16:47:09"

Write-Host "Final used range:" $ws.UsedRange.Address()
